# NOTE: worksheet object references captured before a later Worksheets.Add()/
# Move() call become stale in this host, so every sheet handle is re-fetched
# via Worksheets.Item(name) immediately before it is used.

$wb = $excel.ActiveWorkbook

# --- Sheet1: replace the original text content with numeric data ---
$wb.Worksheets.Item("Sheet1").Range("A1").Value = 3424
$wb.Worksheets.Item("Sheet1").Range("B1").Value = 23423
$wb.Worksheets.Item("Sheet1").Range("C1").Value = 24234
$wb.Worksheets.Item("Sheet1").Range("A2").Value = 234234
$wb.Worksheets.Item("Sheet1").Range("B2").Value = 34234
$wb.Worksheets.Item("Sheet1").Range("C2").Value = 23424

$wb.Worksheets.Item("Sheet1").Range("C6").Select()

# --- Add Sheet2 right after Sheet1 ---
$null = $wb.Worksheets.Add()
$wb.Worksheets.Item(1).Name = "Sheet2"
$wb.Worksheets.Item("Sheet2").Move($null, $wb.Worksheets.Item("Sheet1"))

$wb.Worksheets.Item("Sheet2").Range("A1").Value = "aksdj"
$wb.Worksheets.Item("Sheet2").Range("B1").Value = "sjdf"
$wb.Worksheets.Item("Sheet2").Range("C1").Value = "98usdf"
$wb.Worksheets.Item("Sheet2").Range("A2").Value = "lksj"
$wb.Worksheets.Item("Sheet2").Range("B2").Value = "iu89sd"
$wb.Worksheets.Item("Sheet2").Range("C2").Value = "879dfk"

$wb.Worksheets.Item("Sheet2").Range("D12").Select()

# --- Add Sheet3 right after Sheet2 ---
$null = $wb.Worksheets.Add()
$wb.Worksheets.Item(1).Name = "Sheet3"
$wb.Worksheets.Item("Sheet3").Move($null, $wb.Worksheets.Item("Sheet2"))

$wb.Worksheets.Item("Sheet3").Range("A1").Value = "sdfsdf23"
$wb.Worksheets.Item("Sheet3").Range("B1").Value = "sdfsdf"
$wb.Worksheets.Item("Sheet3").Range("C1").Value = "xv23r"
$wb.Worksheets.Item("Sheet3").Range("A2").Value = "gvdfg"
$wb.Worksheets.Item("Sheet3").Range("B2").Value = "sdfgsgs"
$wb.Worksheets.Item("Sheet3").Range("C2").Value = "sfsdf"

$wb.Worksheets.Item("Sheet3").Range("C2").Select()

# --- Sheet2 ends up the active/selected tab ---
$wb.Worksheets.Item("Sheet2").Activate()
